# Auto-generated Excel COM-interop script applying numeric corrections
# to the market-price / profit columns (H,I,J,K,L,M,N) across several
# worksheets, per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 36770
$ws.Range("J3").Value = 36770
$ws.Range("L3").Value = 36770
$ws.Range("N3").Value = -36998
# Row 17
$ws.Range("H17").Value = 1881.5834
$ws.Range("J17").Value = 1881.5834
$ws.Range("L17").Value = 5644.7502
$ws.Range("N17").Value = -5980.7502
# Row 31
$ws.Range("H31").Value = 8000
$ws.Range("I31").Value = 8000
$ws.Range("K31").Value = 24000
$ws.Range("M31").Value = -23770
# Row 40
$ws.Range("H40").Value = 1399.3334
$ws.Range("I40").Value = 1349.5
$ws.Range("J40").Value = 1499
$ws.Range("K40").Value = 1349.5
$ws.Range("L40").Value = 1499
$ws.Range("M40").Value = -1174.5
$ws.Range("N40").Value = -1849
# Row 44
$ws.Range("H44").Value = 140001
$ws.Range("J44").Value = 140001
$ws.Range("L44").Value = 140001
$ws.Range("N44").Value = -140925
# Row 102
$ws.Range("H102").Value = 36770
$ws.Range("J102").Value = 36770
$ws.Range("L102").Value = 36770
$ws.Range("N102").Value = -43260
# Row 110
$ws.Range("H110").Value = 94999.5
$ws.Range("J110").Value = 94999.5
$ws.Range("L110").Value = 94999.5
$ws.Range("N110").Value = -103179.5
# Row 132
$ws.Range("H132").Value = 85169.336
$ws.Range("I132").Value = 101603.3
$ws.Range("K132").Value = 304809.9
$ws.Range("M132").Value = -302279.9
# Row 137
$ws.Range("H137").Value = 2138.5715
$ws.Range("I137").Value = 1924.25
$ws.Range("K137").Value = 5772.75
$ws.Range("M137").Value = -3222.75
# Row 138
$ws.Range("H138").Value = 4410.276
$ws.Range("I138").Value = 1378.8
$ws.Range("K138").Value = 4136.4
$ws.Range("M138").Value = 1003.6

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2832.484
$ws.Range("I32").Value = 2993.2144
$ws.Range("K32").Value = 2993.2144
$ws.Range("M32").Value = -2706.2144
# Row 63
$ws.Range("H63").Value = 3697.25
$ws.Range("I63").Value = 3999
$ws.Range("K63").Value = 3999
$ws.Range("M63").Value = -3313
# Row 66
$ws.Range("H66").Value = 3697.25
$ws.Range("I66").Value = 3999
$ws.Range("K66").Value = 19995
$ws.Range("M66").Value = -16563
# Row 74
$ws.Range("H74").Value = 5790399
$ws.Range("I74").Value = 3090402.8
$ws.Range("J74").Value = 13890388
$ws.Range("K74").Value = 3090402.8
$ws.Range("L74").Value = 13890388
$ws.Range("M74").Value = -3089528.8
$ws.Range("N74").Value = -13892136
# Row 77
$ws.Range("H77").Value = 5790399
$ws.Range("I77").Value = 3090402.8
$ws.Range("J77").Value = 13890388
$ws.Range("K77").Value = 15452014
$ws.Range("L77").Value = 69451940
$ws.Range("M77").Value = -15447646
$ws.Range("N77").Value = -69460676
# Row 110
$ws.Range("H110").Value = 1536.6364
$ws.Range("I110").Value = 1536.6364
$ws.Range("K110").Value = 1536.6364
$ws.Range("M110").Value = 508.3635999999999
# Row 132
$ws.Range("H132").Value = 2600.4722
$ws.Range("I132").Value = 2625.8667
$ws.Range("J132").Value = 2473.5
$ws.Range("K132").Value = 7877.6001
$ws.Range("L132").Value = 7420.5
$ws.Range("M132").Value = -5347.6001
$ws.Range("N132").Value = -12480.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4049.3
$ws.Range("I105").Value = 2149.25
$ws.Range("J105").Value = 5316
$ws.Range("K105").Value = 2149.25
$ws.Range("L105").Value = 5316
$ws.Range("M105").Value = -402.25
$ws.Range("N105").Value = -8810
# Row 134
$ws.Range("H134").Value = 18231898
$ws.Range("I134").Value = 8930900
$ws.Range("K134").Value = 26792700
$ws.Range("M134").Value = -26790165

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1117.8
$ws.Range("J16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -2074
# Row 17
$ws.Range("H17").Value = 9672.333000000001
$ws.Range("J17").Value = 9672.333000000001
$ws.Range("L17").Value = 9672.333000000001
$ws.Range("N17").Value = -10020.333
# Row 25
$ws.Range("H25").Value = 7506.5
$ws.Range("J25").Value = 7506.5
$ws.Range("L25").Value = 7506.5
$ws.Range("N25").Value = -7854.5
# Row 31
$ws.Range("H31").Value = 2971.158
$ws.Range("I31").Value = 2937.4546
$ws.Range("J31").Value = 3017.5
$ws.Range("K31").Value = 2937.4546
$ws.Range("L31").Value = 3017.5
$ws.Range("M31").Value = -2642.4546
$ws.Range("N31").Value = -3607.5
# Row 34
$ws.Range("H34").Value = 2971.158
$ws.Range("I34").Value = 2937.4546
$ws.Range("J34").Value = 3017.5
$ws.Range("K34").Value = 2937.4546
$ws.Range("L34").Value = 3017.5
$ws.Range("M34").Value = -2735.4546
$ws.Range("N34").Value = -3421.5
# Row 58
$ws.Range("H58").Value = 2163.3333
$ws.Range("I58").Value = 2084.8696
$ws.Range("K58").Value = 2084.8696
$ws.Range("M58").Value = -1881.8696
# Row 59
$ws.Range("H59").Value = 37831.94
$ws.Range("J59").Value = 38633.938
$ws.Range("L59").Value = 38633.938
$ws.Range("N59").Value = -40923.938
# Row 113
$ws.Range("H113").Value = 1117.8
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840
# Row 122
$ws.Range("H122").Value = 2785.75
$ws.Range("J122").Value = 2774.5
$ws.Range("L122").Value = 8323.5
$ws.Range("N122").Value = -13223.5
# Row 132
$ws.Range("H132").Value = 11262.593
$ws.Range("I132").Value = 11443.96
$ws.Range("K132").Value = 34331.88
$ws.Range("M132").Value = -31801.88
# Row 134
$ws.Range("H134").Value = 4548118.5
$ws.Range("I134").Value = 2579
$ws.Range("J134").Value = 16669557
$ws.Range("K134").Value = 7737
$ws.Range("L134").Value = 50008671
$ws.Range("M134").Value = -5202
$ws.Range("N134").Value = -50013741
# Row 136
$ws.Range("H136").Value = 2163.3333
$ws.Range("I136").Value = 2084.8696
$ws.Range("K136").Value = 6254.6088
$ws.Range("M136").Value = -3704.6088

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 204240
$ws.Range("I4").Value = 5300
$ws.Range("K4").Value = 15900
$ws.Range("M4").Value = -15788
# Row 11
$ws.Range("H11").Value = 418.8889
$ws.Range("I11").Value = 358.875
$ws.Range("K11").Value = 1076.625
$ws.Range("M11").Value = -936.625
# Row 59
$ws.Range("H59").Value = 849
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# Row 121
$ws.Range("H121").Value = 2361136.2
$ws.Range("J121").Value = 4346550
$ws.Range("L121").Value = 13039650
$ws.Range("N121").Value = -13042270
# Row 137
$ws.Range("H137").Value = 281495.97
$ws.Range("J137").Value = 378819.8
$ws.Range("L137").Value = 1136459.4
$ws.Range("N137").Value = -1146659.4

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 9537.25
$ws.Range("I70").Value = 9419.049999999999
$ws.Range("J70").Value = 9685
$ws.Range("K70").Value = 9419.049999999999
$ws.Range("L70").Value = 9685
$ws.Range("M70").Value = -9149.049999999999
$ws.Range("N70").Value = -10225
# Row 73
$ws.Range("H73").Value = 9537.25
$ws.Range("I73").Value = 9419.049999999999
$ws.Range("J73").Value = 9685
$ws.Range("K73").Value = 9419.049999999999
$ws.Range("L73").Value = 9685
$ws.Range("M73").Value = -8483.049999999999
$ws.Range("N73").Value = -11557
# Row 122
$ws.Range("H122").Value = 2807.3914
$ws.Range("I122").Value = 3175.9333
$ws.Range("K122").Value = 9527.7999
$ws.Range("M122").Value = -7077.7999
# Row 126
$ws.Range("H126").Value = 8709.615
$ws.Range("I126").Value = 8835.75
$ws.Range("K126").Value = 26507.25
$ws.Range("M126").Value = -24037.25
# Row 132
$ws.Range("H132").Value = 3356.6667
$ws.Range("I132").Value = 3401.25
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 10203.75
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -7673.75
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 6.9166665
$ws.Range("I12").Value = 6.9166665
$ws.Range("K12").Value = 6.9166665
$ws.Range("M12").Value = 163.0833335
# Row 45
$ws.Range("H45").Value = 12000
$ws.Range("I45").Value = 10000
$ws.Range("K45").Value = 10000
$ws.Range("M45").Value = -9593
# Row 136
$ws.Range("H136").Value = 22730270
$ws.Range("I136").Value = 2560.8718
$ws.Range("J136").Value = 200006400
$ws.Range("K136").Value = 7682.6154
$ws.Range("L136").Value = 600019200
$ws.Range("M136").Value = -5132.6154
$ws.Range("N136").Value = -600024300

$ws = $wb.Worksheets.Item("WVR")
# Row 31
$ws.Range("H31").Value = 17008
$ws.Range("I31").Value = 4017
$ws.Range("J31").Value = 29999
$ws.Range("K31").Value = 4017
$ws.Range("L31").Value = 29999
$ws.Range("M31").Value = -3669
$ws.Range("N31").Value = -30695
# Row 81
$ws.Range("H81").Value = 2081.5
$ws.Range("I81").Value = 2199.6667
$ws.Range("K81").Value = 4399.3334
$ws.Range("M81").Value = -3338.3334
# Row 84
$ws.Range("H84").Value = 2081.5
$ws.Range("I84").Value = 2199.6667
$ws.Range("K84").Value = 21996.667
$ws.Range("M84").Value = -16692.667
# Row 126
$ws.Range("H126").Value = 1233.1666
$ws.Range("I126").Value = 879.8
$ws.Range("K126").Value = 2639.4
$ws.Range("M126").Value = -169.3999999999996
# Row 132
$ws.Range("H132").Value = 334.83334
$ws.Range("I132").Value = 244.25
$ws.Range("J132").Value = 516
$ws.Range("K132").Value = 732.75
$ws.Range("L132").Value = 1548
$ws.Range("M132").Value = 1797.25
$ws.Range("N132").Value = -6608
# Row 136
$ws.Range("H136").Value = 1601.7059
$ws.Range("I136").Value = 1373.5714
$ws.Range("J136").Value = 2666.3333
$ws.Range("K136").Value = 4120.7142
$ws.Range("L136").Value = 7998.999899999999
$ws.Range("M136").Value = -1570.7142
$ws.Range("N136").Value = -13098.9999

